$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.923.23'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.554.13'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.72'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '1.554.49'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.73'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '26.920.33'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.93'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.28'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.18'
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.07'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.91'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = '1.435.86'
$ws.Range("E33").Value = '  +4.88%  '
$ws.Range("E34").Value = '  +3.98%  '
$ws.Range("E35").Value = '  +3.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.810'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.986'
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.98'
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("D47").Value = '1.689.83'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.17'
$ws.Range("E48").Value = '  +3.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0524'
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +3.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0955'
$ws.Range("E51").Value = '  +1.59%  '
